# This edit reorders the data rows (rows 2-24) of the "Artfynd" sheet.
# The row order produced by an upstream export changed, so each data
# row's full contents (columns A:AY) need to move to a new row position.
# Mapping below: key = destination row (after), value = source row (before),
# expressed against the original (pre-edit) row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 24
$firstCol = 1   # A
$lastCol = 51   # AY

# destinationRow -> sourceRow (both referring to the ORIGINAL layout)
$mapping = @{
    2  = 6
    3  = 3
    4  = 5
    5  = 2
    6  = 4
    7  = 10
    8  = 17
    9  = 20
    10 = 7
    11 = 19
    12 = 22
    13 = 11
    14 = 21
    15 = 24
    16 = 18
    17 = 9
    18 = 8
    19 = 13
    20 = 23
    21 = 15
    22 = 14
    23 = 16
    24 = 12
}

# Snapshot every cell value of the original data rows before overwriting anything,
# since several destination rows read from rows that will themselves be overwritten.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Columns Y (Startdatum) and AA (Slutdatum) hold plain text dates
# ("yyyy-mm-dd"). Force text formatting on them first so assigning the
# string back doesn't get auto-coerced into a real Excel date serial.
$ws.Range("Y2:Y24").NumberFormat = "@"
$ws.Range("AA2:AA24").NumberFormat = "@"

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
    }
}
